$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Km initiali: value update
$ws.Range("B12").Value = 256004

# Row 16 (Ziua 3)
$ws.Range("B16").Value = 30
$ws.Range("C16").Value = "Acasa-Birou"
$ws.Range("D16").Value = " "

# Row 19 (Ziua 6)
$ws.Range("B19").Value = 30
$ws.Range("C19").Value = "Acasa-Birou"
$ws.Range("D19").Value = " "

# Row 20 (Ziua 7)
$ws.Range("B20").Value = 85
$ws.Range("C20").Value = "Cluj-Apahida"
$ws.Range("D20").Value = "Interes Serviciu"

# Row 21 (Ziua 8)
$ws.Range("B21").Value = 101
$ws.Range("C21").Value = "Cluj-Dej"
$ws.Range("D21").Value = "Interes Serviciu"

# Row 22 (Ziua 9)
$ws.Range("B22").Value = 257
$ws.Range("C22").Value = "Cluj-Bistrita"
$ws.Range("D22").Value = "Interes Serviciu"

# Row 23 (Ziua 10)
$ws.Range("B23").Value = 421
$ws.Range("C23").Value = "Cluj-Satu-Mare"
$ws.Range("D23").Value = "Interes Serviciu"

# Row 27 (Ziua 14)
$ws.Range("B27").Value = 30
$ws.Range("C27").Value = "Acasa-Birou"
$ws.Range("D27").Value = " "

# Row 28 (Ziua 15)
$ws.Range("B28").Value = 121
$ws.Range("C28").Value = "Cluj-Turda"
$ws.Range("D28").Value = "Interes Serviciu"

# Row 29 (Ziua 16)
$ws.Range("B29").Value = 421
$ws.Range("C29").Value = "Cluj-Satu-Mare"
$ws.Range("D29").Value = "Interes Serviciu"

# Row 30 (Ziua 17)
$ws.Range("B30").Value = 47
$ws.Range("C30").Value = "Cluj-Cluj"
$ws.Range("D30").Value = "Interes Serviciu"

# Row 33 (Ziua 20)
$ws.Range("B33").Value = 30
$ws.Range("C33").Value = "Acasa-Birou"
$ws.Range("D33").Value = " "

# Row 34 (Ziua 21)
$ws.Range("B34").Value = 356
$ws.Range("C34").Value = "Cluj-Baia-Mare"
$ws.Range("D34").Value = "Interes Serviciu"

# Row 36 (Ziua 23)
$ws.Range("B36").Value = 30
$ws.Range("C36").Value = "Acasa-Birou"
$ws.Range("D36").Value = " "

# Row 37 (Ziua 24)
$ws.Range("B37").Value = 92
$ws.Range("C37").Value = "Cluj-Bontida"
$ws.Range("D37").Value = "Interes Serviciu"

# Row 40 (Ziua 27)
$ws.Range("B40").Value = 156
$ws.Range("C40").Value = "Cluj-Zalau"
$ws.Range("D40").Value = "Interes Serviciu"

# Row 41 (Ziua 28)
$ws.Range("B41").Value = 85
$ws.Range("C41").Value = "Cluj-Apahida"
$ws.Range("D41").Value = "Interes Serviciu"

# Row 42 (Ziua 29)
$ws.Range("B42").Value = 30
$ws.Range("C42").Value = "Acasa-Birou"
$ws.Range("D42").Value = " "

# Row 43 (Ziua 30)
$ws.Range("B43").Value = 152
$ws.Range("C43").Value = "Cluj-Cmp. Turzii"
$ws.Range("D43").Value = "Interes Serviciu"

# Totals
$ws.Range("B44").Value = 2534
$ws.Range("B45").Value = 258538
